$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 02:21"

# Swap country labels for rows 157/158 (Burkina Faso overtakes Uruguay in ranking)
$ws.Range("A157").Value = "Burkina Faso"
$ws.Range("A158").Value = "Uruguay"

# Update per-country statistics (columns B,C,D,E,F,G,H)
# Row 4
$ws.Range("B4").Value = 7943598
$ws.Range("C4").Value = 48969
$ws.Range("D4").Value = 5085449
$ws.Range("E4").Value = 2638895
$ws.Range("G4").Value = 607
$ws.Range("H4").Value = 219254
# Row 6
$ws.Range("B6").Value = 5091840
$ws.Range("C6").Value = 34650
$ws.Range("E6").Value = 487882
$ws.Range("G6").Value = 544
$ws.Range("H6").Value = 150236
# Row 25
$ws.Range("B25").Value = 323453
$ws.Range("C25").Value = 2975
$ws.Range("E25").Value = 40262
# Row 38
$ws.Range("B38").Value = 119666
$ws.Range("C38").Value = 825
$ws.Range("D38").Value = 95552
$ws.Range("E38").Value = 21632
$ws.Range("G38").Value = 8
$ws.Range("H38").Value = 2482
# Row 40
$ws.Range("B40").Value = 114005
$ws.Range("C40").Value = 4631
$ws.Range("D40").Value = 54909
$ws.Range("E40").Value = 58148
$ws.Range("G40").Value = 43
$ws.Range("H40").Value = 948
# Row 111
$ws.Range("B111").Value = 9514
$ws.Range("C111").Value = 154
$ws.Range("E111").Value = 1345
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 131
# Row 131
$ws.Range("B131").Value = 5035
$ws.Range("C131").Value = 17
$ws.Range("D131").Value = 4830
$ws.Range("E131").Value = 98
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 107
# Row 134
$ws.Range("B134").Value = 4854
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 1924
$ws.Range("E134").Value = 2868
# Row 157
$ws.Range("B157").Value = 2271
$ws.Range("D157").Value = 1542
$ws.Range("E157").Value = 668
$ws.Range("H157").Value = 61
# Row 158
$ws.Range("B158").Value = 2268
$ws.Range("C158").Value = 17
$ws.Range("D158").Value = 1930
$ws.Range("E158").Value = 288
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 50
# Row 159
$ws.Range("E159").Value = 127
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 595
# Row 160
$ws.Range("B160").Value = 1986
$ws.Range("C160").Value = 35
$ws.Range("D160").Value = 1444
$ws.Range("E160").Value = 517
# Row 173
$ws.Range("B173").Value = 696
$ws.Range("C173").Value = 1
$ws.Range("E173").Value = 18
# Row 190
$ws.Range("B190").Value = 221
$ws.Range("C190").Value = 1
$ws.Range("E190").Value = 9
